# NATMI output update (Fgf2-Sdc4), following Dr Hou advice:
# the "ECs" sending cluster is added alongside the existing "FAPs" and
# "sCs" clusters, and every derived NATMI metric (columns E:T) is
# recomputed for the resulting 3x3 sending/target cluster combinations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row below is: row number, then the values for columns A..T, in order
$rowsData = @(
    @(2, "ECs", "Fgf2", "Sdc4", "ECs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 3, 1, 2.300909333333333, 6.902728, 0.03776979643482627, 0.03776979643482627, 1.486938113633778, 13.382443022704, 0.001307533900438009, 0.001307533900438008),
    @(3, "ECs", "Fgf2", "Sdc4", "FAPs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 3, 1, 16.96312166666667, 50.889365, 0.2784523679257784, 0.2784523679257784, 10.96223643711889, 98.66012793406999, 0.009639604792375344, 0.009639604792375342),
    @(4, "ECs", "Fgf2", "Sdc4", "sCs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 3, 1, 41.655263, 124.965789, 0.6837778356393953, 0.6837778356393953, 26.91926939094466, 242.273424518502, 0.02367136667017492, 0.02367136667017491),
    @(5, "FAPs", "Fgf2", "Sdc4", "ECs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 3, 1, 2.300909333333333, 6.902728, 0.03776979643482627, 0.03776979643482627, 35.25643335644267, 317.307900207984, 0.03100262304086453, 0.03100262304086452),
    @(6, "FAPs", "Fgf2", "Sdc4", "FAPs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 3, 1, 16.96312166666667, 50.889365, 0.2784523679257784, 0.2784523679257784, 259.9229617151634, 2339.30665543647, 0.2285623596763431, 0.2285623596763431),
    @(7, "FAPs", "Fgf2", "Sdc4", "sCs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 3, 1, 41.655263, 124.965789, 0.6837778356393953, 0.6837778356393953, 638.276346933238, 5744.487122399142, 0.5612661036870082, 0.5612661036870081),
    @(8, "sCs", "Fgf2", "Sdc4", "ECs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 3, 1, 2.300909333333333, 6.902728, 0.03776979643482627, 0.03776979643482627, 6.208746134154667, 55.878715207392, 0.005459639493523736, 0.005459639493523735),
    @(9, "sCs", "Fgf2", "Sdc4", "FAPs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 3, 1, 16.96312166666667, 50.889365, 0.2784523679257784, 0.2784523679257784, 45.77308394787333, 411.95775553086, 0.04025040345705995, 0.04025040345705995),
    @(10, "sCs", "Fgf2", "Sdc4", "sCs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 3, 1, 41.655263, 124.965789, 0.6837778356393953, 0.6837778356393953, 112.402061816044, 1011.618556344396, 0.0988403652822122, 0.0988403652822122)
)

foreach ($entry in $rowsData) {
    $rowNum = $entry[0]
    for ($i = 1; $i -lt $entry.Length; $i++) {
        $ws.Cells.Item($rowNum, $i).Value = $entry[$i]
    }
}

# Grow the declared sheet dimension to cover the newly added rows
$ws.Range("A1:T10").Select() | Out-Null